$wb = $excel.ActiveWorkbook

# --- Update the "Logs" worksheet: append a new row (row 11) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Wil je dit oppakken?"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Testmail #2: Wil je dit oppakken?"
$logs.Range("D11").Value = "Overig"
$logs.Range("E11").Value = "Geachte afzender,`nDank u voor uw e-mail. Kunt u wat meer context geven over wat u precies bedoelt met ""Testmail #2: Wil je dit oppakken""? Zo kan ik u beter van dienst zijn. `nMet vriendelijke groet,`n[Naam]`nNederlandse e-mailassistent"
$logs.Range("F11").Value = "2025-08-03 14:28:55"
$logs.Range("G11").Value = "Ja"
$logs.Range("H11").Value = "Nee"
$logs.Range("I11").Value = "Ja"
$logs.Range("J11").Value = "Nee"
$logs.Rows.Item(11).AutoFit()

# --- Extend the conditional formatting ranges to include the new row ---
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2:" + $col + "10")
    $newRange = $logs.Range($col + "2:" + $col + "11")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# --- Update the "Dashboard" worksheet: Overig count goes from 1 to 2 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B4").Value = 2
